$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.805635333333334
$ws.Range("H2").Value = 5.416906000000001
$ws.Range("I2").Value = 0.01900969238460649
$ws.Range("J2").Value = 0.01900969238460649
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.083576666666666
$ws.Range("N2").Value = 9.250729999999999
$ws.Range("O2").Value = 0.2272509363535097
$ws.Range("P2").Value = 0.2272509363535097
$ws.Range("Q2").Value = 5.567814982375555
$ws.Range("R2").Value = 50.11033484138
$ws.Range("S2").Value = 0.004319970394194007
$ws.Range("T2").Value = 0.004319970394194007
$ws.Range("G3").Value = 1.805635333333334
$ws.Range("H3").Value = 5.416906000000001
$ws.Range("I3").Value = 0.01900969238460649
$ws.Range("J3").Value = 0.01900969238460649
$ws.Range("M3").Value = 6.453984666666667
$ws.Range("O3").Value = 0.4756405360586227
$ws.Range("P3").Value = 0.4756405360586227
$ws.Range("Q3").Value = 11.65354275492489
$ws.Range("R3").Value = 104.881884794324
$ws.Range("S3").Value = 0.00904178027612375
$ws.Range("T3").Value = 0.009041780276123748
$ws.Range("G4").Value = 1.805635333333334
$ws.Range("H4").Value = 5.416906000000001
$ws.Range("I4").Value = 0.01900969238460649
$ws.Range("J4").Value = 0.01900969238460649
$ws.Range("M4").Value = 4.031477000000001
$ws.Range("N4").Value = 12.094431
$ws.Range("O4").Value = 0.2971085275878677
$ws.Range("P4").Value = 0.2971085275878677
$ws.Range("Q4").Value = 7.279377316720669
$ws.Range("R4").Value = 65.51439585048603
$ws.Range("S4").Value = 0.005647941714288735
$ws.Range("T4").Value = 0.005647941714288735
$ws.Range("I5").Value = 0.7995527524661065
$ws.Range("J5").Value = 0.7995527524661064
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.083576666666666
$ws.Range("N5").Value = 9.250729999999999
$ws.Range("O5").Value = 0.2272509363535097
$ws.Range("P5").Value = 0.2272509363535097
$ws.Range("Q5").Value = 234.1837892119344
$ws.Range("R5").Value = 2107.65410290741
$ws.Range("S5").Value = 0.1816991116619486
$ws.Range("T5").Value = 0.1816991116619486
$ws.Range("I6").Value = 0.7995527524661065
$ws.Range("J6").Value = 0.7995527524661064
$ws.Range("M6").Value = 6.453984666666667
$ws.Range("O6").Value = 0.4756405360586227
$ws.Range("P6").Value = 0.4756405360586227
$ws.Range("Q6").Value = 490.1511290749132
$ws.Range("S6").Value = 0.3802996997901262
$ws.Range("T6").Value = 0.3802996997901261
$ws.Range("I7").Value = 0.7995527524661065
$ws.Range("J7").Value = 0.7995527524661064
$ws.Range("M7").Value = 4.031477000000001
$ws.Range("N7").Value = 12.094431
$ws.Range("O7").Value = 0.2971085275878677
$ws.Range("P7").Value = 0.2971085275878677
$ws.Range("Q7").Value = 306.1725593485364
$ws.Range("R7").Value = 2755.553034136828
$ws.Range("S7").Value = 0.2375539410140317
$ws.Range("T7").Value = 0.2375539410140317
$ws.Range("G8").Value = 17.23384333333334
$ws.Range("H8").Value = 51.70153000000001
$ws.Range("I8").Value = 0.1814375551492871
$ws.Range("J8").Value = 0.1814375551492871
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.083576666666666
$ws.Range("N8").Value = 9.250729999999999
$ws.Range("O8").Value = 0.2272509363535097
$ws.Range("P8").Value = 0.2272509363535097
$ws.Range("Q8").Value = 53.14187717965556
$ws.Range("R8").Value = 478.2768946169
$ws.Range("S8").Value = 0.04123185429736703
$ws.Range("T8").Value = 0.04123185429736703
$ws.Range("G9").Value = 17.23384333333334
$ws.Range("H9").Value = 51.70153000000001
$ws.Range("I9").Value = 0.1814375551492871
$ws.Range("J9").Value = 0.1814375551492871
$ws.Range("M9").Value = 6.453984666666667
$ws.Range("O9").Value = 0.4756405360586227
$ws.Range("P9").Value = 0.4756405360586227
$ws.Range("Q9").Value = 111.2269606210689
$ws.Range("R9").Value = 1001.04264558962
$ws.Range("S9").Value = 0.08629905599237282
$ws.Range("T9").Value = 0.08629905599237281
$ws.Range("G10").Value = 17.23384333333334
$ws.Range("H10").Value = 51.70153000000001
$ws.Range("I10").Value = 0.1814375551492871
$ws.Range("J10").Value = 0.1814375551492871
$ws.Range("M10").Value = 4.031477000000001
$ws.Range("N10").Value = 12.094431
$ws.Range("O10").Value = 0.2971085275878677
$ws.Range("P10").Value = 0.2971085275878677
$ws.Range("Q10").Value = 69.47784301993669
$ws.Range("R10").Value = 625.3005871794302
$ws.Range("S10").Value = 0.05390664485954721
$ws.Range("T10").Value = 0.05390664485954721
